$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.018.32"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.333.18"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.35"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.14"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.25"
$ws.Range("E10").Value = "  -3.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.27"
$ws.Range("E11").Value = "  +2.47%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("E13").Value = "  +3.48%  "
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").Value = "2.693.46"
$ws.Range("D16").Value = "2.309.93"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.793"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "42.964.49"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.16"
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.07"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.98"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.71"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("E28").Value = "  -13.74%  "
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.65"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "138.51"
$ws.Range("E32").Value = "  -16.50%  "
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.73"
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.40"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("E37").Value = "  +3.19%  "
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.24"
$ws.Range("E41").Value = "  +23.69%  "
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "1.934.88"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.05"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.89"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").Value = "2.561.58"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.66"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.04"
$ws.Range("E51").Value = "  +1.77%  "
